$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7566934823989868
$ws.Range("B1").Value = 1.087204217910767
$ws.Range("C1").Value = 2.065640449523926
$ws.Range("D1").Value = 3.516237735748291
$ws.Range("E1").Value = 3.398043870925903
